$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (7): shifts G:J -> H:K
$ws.Columns.Item(7).Insert()

# Set width of the new column G
$ws.Columns.Item(7).ColumnWidth = 22.85546875

# New header cell G1
$ws.Range("G1").Value = "WideStringCol"

# New data cells in column G
$ws.Range("G2").Value = "ABC"
$ws.Range("G3").Value = "äöüß"
$ws.Range("G6").Value = "123adf"
$ws.Range("G7").Value = "dfsf"
$ws.Range("G8").Value = "äöü"

# New cell E3
$ws.Range("E3").Value = "äöüß"

# Update selection to E3
$ws.Range("E3").Select()
